$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell 2 4 "43.105.55"
Set-TextCell 2 5 "  -1.46%  "
Set-TextCell 3 4 "2.274.71"
Set-TextCell 3 5 "  -0.59%  "
Set-TextCell 4 5 "  -0.05%  "
Set-TextCell 5 4 "111.45"
Set-TextCell 5 5 "  +0.40%  "
Set-TextCell 6 4 "263.87"
Set-TextCell 6 5 "  -1.75%  "
Set-TextCell 7 4 "0.648"
Set-TextCell 7 5 "  +3.80%  "
Set-TextCell 8 5 "  -0.25%  "
Set-TextCell 9 4 "0.609"
Set-TextCell 9 5 "  -0.99%  "
Set-TextCell 10 4 "46.58"
Set-TextCell 10 5 "  -2.13%  "
Set-TextCell 11 4 "0.0935"
Set-TextCell 11 5 "  -1.19%  "
Set-TextCell 12 4 "9.12"
Set-TextCell 12 5 "  +1.35%  "
Set-TextCell 13 5 "  +1.18%  "
Set-TextCell 14 4 "15.25"
Set-TextCell 14 5 "  -2.92%  "
Set-TextCell 15 4 "2.619.12"
Set-TextCell 15 5 "  -0.49%  "
Set-TextCell 16 4 "0.863"
Set-TextCell 16 5 "  +2.13%  "
Set-TextCell 17 4 "2.277.23"
Set-TextCell 17 5 "  -0.40%  "
Set-TextCell 18 4 "43.180.23"
Set-TextCell 18 5 "  -1.03%  "
Set-TextCell 19 4 "0.0000109"
Set-TextCell 19 5 "  +0.28%  "
Set-TextCell 20 4 "6.68"
Set-TextCell 20 5 "  -0.48%  "
Set-TextCell 21 4 "72.25"
Set-TextCell 21 5 "  +0.08%  "
Set-TextCell 22 4 "2.40"
Set-TextCell 22 5 "  -1.21%  "
Set-TextCell 23 4 "235.49"
Set-TextCell 23 5 "  +1.61%  "
Set-TextCell 24 5 "  +3.64%  "
Set-TextCell 25 4 "9.32"
Set-TextCell 25 5 "  -5.77%  "
Set-TextCell 26 5 "  +1.89%  "
Set-TextCell 27 5 "  -0.55%  "
Set-TextCell 28 4 "41.03"
Set-TextCell 28 5 "  -1.59%  "
Set-TextCell 29 4 "3.40"
Set-TextCell 29 5 "  +0.12%  "
Set-TextCell 30 4 "2.25"
Set-TextCell 30 5 "  -0.99%  "
Set-TextCell 31 4 "173.43"
Set-TextCell 31 5 "  -1.01%  "
Set-TextCell 32 4 "21.65"
Set-TextCell 32 5 "  +0.95%  "
Set-TextCell 33 5 "  -2.31%  "
Set-TextCell 34 4 "5.59"
Set-TextCell 34 5 "  -0.05%  "
Set-TextCell 35 5 "  +2.55%  "
Set-TextCell 36 4 "0.0376"
Set-TextCell 36 5 "  +4.09%  "
Set-TextCell 37 4 "4.69"
Set-TextCell 37 5 "  +0.45%  "
Set-TextCell 38 4 "3.80"
Set-TextCell 38 5 "  -0.72%  "
Set-TextCell 39 5 "  -2.58%  "
Set-TextCell 40 4 "2.55"
Set-TextCell 40 5 "  +6.52%  "
Set-TextCell 41 4 "14.43"
Set-TextCell 41 5 "  +5.68%  "
Set-TextCell 42 4 "74.03"
Set-TextCell 42 5 "  +0.63%  "
Set-TextCell 43 5 "  -2.93%  "
Set-TextCell 44 2 "FirstDigitalUSD"
Set-TextCell 44 3 "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell 44 4 "1.00"
Set-TextCell 44 5 "  +0.00%  "
Set-TextCell 45 2 "THORChain"
Set-TextCell 45 3 "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextCell 45 4 "5.94"
Set-TextCell 45 5 "  -5.28%  "
Set-TextCell 46 5 "  -1.73%  "
Set-TextCell 47 4 "8.55"
Set-TextCell 47 5 "  -1.81%  "
Set-TextCell 48 5 "  +3.46%  "
Set-TextCell 49 4 "0.0998"
Set-TextCell 49 5 "  +0.43%  "
Set-TextCell 50 2 "Aave"
Set-TextCell 50 3 "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell 50 4 "100.23"
Set-TextCell 50 5 "  -1.98%  "
Set-TextCell 51 2 "ordi"
Set-TextCell 51 3 "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
Set-TextCell 51 4 "70.31"
Set-TextCell 51 5 "  +28.76%  "
